$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.159.82"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.559.86"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'206.44"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'0.491"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'22.20"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'0.0592"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "'0.0861"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.781.90"
$ws.Range("D13").Value = "1.558.72"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "'3.77"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "'0.515"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("D16").Value = "'62.80"
$ws.Range("E16").Value = "  -1.03%  "
$ws.Range("D17").Value = "27.161.39"
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "'214.00"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.11"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").Value = "'9.35"
$ws.Range("E23").Value = "  -3.65%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'152.14"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").Value = "'6.60"
$ws.Range("E26").Value = "  -3.25%  "
$ws.Range("D27").Value = "'14.88"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "'0.0462"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "'3.16"
$ws.Range("E32").Value = "  -2.18%  "
$ws.Range("D33").Value = "1.384.47"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("D39").Value = "'0.810"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'0.515"
$ws.Range("E40").Value = "  -4.09%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'0.982"
$ws.Range("E42").Value = "  +1.27%  "
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("D44").Value = "'63.29"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").Value = "1.694.96"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'85.45"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").Value = "0.0₇0985"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  +0.03%  "
